# worklog.docx edit: add 2022-04-27 entries, indent the 2022-04-24 entry.
$d = $word.ActiveDocument

function Find-ParaIndex($doc, $substr) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -like "*$substr*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# 1) 2022-04-24 paragraph gets a hanging indent (w:ind w:left="1440"
#    w:hanging="1440" == 72pt left indent / 72pt hanging).
# ---------------------------------------------------------------------
$idx24 = Find-ParaIndex $d "2022-04-24"
if ($idx24 -gt 0) {
    $p24 = $d.Paragraphs.Item($idx24)
    $p24.Format.LeftIndent = 72
    $p24.Format.FirstLineIndent = -72
}

# ---------------------------------------------------------------------
# 2) Insert four new paragraphs after the "... Took me about 2 hours of
#    trial to finish" paragraph (still before the two trailing blank
#    paragraphs at the end of the body).
# ---------------------------------------------------------------------
$idxAnchor = Find-ParaIndex $d "Took me about 2 hours of trial to finish"
if ($idxAnchor -le 0) {
    throw "anchor paragraph not found"
}
$anchorPara = $d.Paragraphs.Item($idxAnchor)

# Make room: one fresh paragraph right after the anchor, which the big
# InsertXML call below will expand into the real set of new paragraphs.
$anchorPara.Range.InsertParagraphAfter()
$newHome = $d.Paragraphs.Item($idxAnchor + 1)
$newRange = $newHome.Range

$payload = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr><w:ind w:left="1428" w:hanging="1428"/></w:pPr>
            <w:r><w:t xml:space="preserve">2022-04-27 - </w:t></w:r>
            <w:r><w:tab/><w:t xml:space="preserve">worked on get_image_path function </w:t></w:r>
            <w:r><w:t>used regex to extract the image name from the url.</w:t></w:r>
            <w:r><w:t xml:space="preserve"> It worked!</w:t></w:r>
            <w:r><w:t xml:space="preserve"> Looked back at scripts from last semester to refresh my memory</w:t></w:r>
            <w:r><w:t xml:space="preserve"> on working with regex.</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Worked on the </w:t></w:r>
            <w:r><w:t xml:space="preserve">image_already_in_db function. Had to find a way to get the </w:t></w:r>
            <w:r><w:t xml:space="preserve">image hash into the query, found I needed a </w:t></w:r>
            <w:r><w:t>parameterized</w:t></w:r>
            <w:r><w:t xml:space="preserve"> query. This link </w:t></w:r>
            <w:r><w:t>XHYPERLINKPLACEHOLDERX</w:t></w:r>
            <w:r><w:t xml:space="preserve"> gave me the info I needed. Took me about 1 hour to figure out.</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Finished the save_image_file </w:t></w:r>
            <w:r><w:t xml:space="preserve">function. This one was pretty simple took me about 20 min. </w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>
            </w:pPr>
            <w:r><w:t>Created the set_desktop_background_image function was able to copy the code from lab 9</w:t></w:r>
            <w:r><w:t>.</w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newRange.InsertXML($payload)

# ---------------------------------------------------------------------
# 3) Turn the placeholder text into a real hyperlink pointing at the
#    stackoverflow question, styled with the built-in Hyperlink style.
# ---------------------------------------------------------------------
$idxLink = Find-ParaIndex $d "XHYPERLINKPLACEHOLDERX"
if ($idxLink -gt 0) {
    $linkPara = $d.Paragraphs.Item($idxLink)
    $find3 = $linkPara.Range.Find
    $find3.ClearFormatting()
    $found3 = $find3.Execute("XHYPERLINKPLACEHOLDERX", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found3) {
        $linkRange = $d.Range($find3.Parent.Start, $find3.Parent.End)
        [void]$d.Hyperlinks.Add($linkRange, "https://stackoverflow.com/questions/45343175/python-3-sqlite-parameterized-sql-query", "", "", "https://stackoverflow.com/questions/45343175/python-3-sqlite-parameterized-sql-query")
    }
}
